$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 44, shifting existing rows 44-173 down to 45-174.
$ws.Rows.Item(44).Insert()

# Populate the newly inserted row 44 with the new weekly data point.
$ws.Range("A44").Value = 8
$ws.Range("B44").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C44").Value = 'Coquimbo'
$ws.Range("D44").Value = 44648
$ws.Range("E44").Value = 4
$ws.Range("F44").Value = 100112037
$ws.Range("G44").Value = 'Cebollín'
$ws.Range("H44").Value = 'Sin especificar'
$ws.Range("I44").Value = 'Primera'
$ws.Range("J44").Value = 2000
$ws.Range("K44").Value = 1100
$ws.Range("L44").Value = 1200
$ws.Range("M44").Value = 1150
$ws.Range("N44").Value = '$/paquete 6 unidades'
$ws.Range("O44").Value = 'Provincia del Elquí'
$ws.Range("P44").Value = 192
$ws.Range("Q44").Value = 6
$ws.Range("R44").Value = 'Hortaliza'
